$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.858.51'
$ws.Range("E2").Value = '  -0.32%  '
$ws.Range("D3").Value = '2.581.61'
$ws.Range("E3").Value = '  +1.26%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.75'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.58'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.596'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +2.43%  '
$ws.Range("E9").Value = '  +2.46%  '
$ws.Range("E10").Value = '  +2.89%  '
$ws.Range("E11").Value = '  -0.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.353'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.24'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.14%  '
$ws.Range("D14").Value = '3.043.97'
$ws.Range("E14").Value = '  +1.27%  '
$ws.Range("D15").Value = '62.780.63'
$ws.Range("E15").Value = '  -0.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000146'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.67%  '
$ws.Range("D17").Value = '2.577.92'
$ws.Range("E17").Value = '  +1.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.31'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '342.03'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.88%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.38'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.67'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.66'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.10'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.74%  '
$ws.Range("D25").Value = '2.713.41'
$ws.Range("E25").Value = '  +1.04%  '
$ws.Range("E26").Value = '  -2.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.59'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.85'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +7.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.31'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.42%  '
$ws.Range("E31").Value = '  -3.04%  '
$ws.Range("E32").Value = '  +2.23%  '
$ws.Range("D33").Value = '0.0₃0819'
$ws.Range("E33").Value = '  +0.91%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '463.36'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +13.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '174.95'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.60'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +3.75%  '
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.400'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.01'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.51'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.99%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.70'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '158.22'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +4.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.76'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.637'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +5.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.12'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0540'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0966'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0236'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.73%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.43'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.71'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.41%  '
